$d = $word.ActiveDocument

$replacements = @(
    @("970×4=3880", "825×8=6600"),
    @("493×7=3451", "559×8=4472"),
    @("171×2=342",  "489×7=3423"),
    @("367×3=1101", "364×7=2548"),
    @("572×7=4004", "187×8=1496"),
    @("987×4=3948", "575×6=3450"),
    @("981×3=2943", "583×5=2915"),
    @("117×9=1053", "273×3=819"),
    @("588×5=2940", "777×3=2331"),
    @("218×5=1090", "539×4=2156"),
    @("471×7=3297", "870×2=1740"),
    @("491×8=3928", "875×9=7875"),
    @("548×3=1644", "382×3=1146"),
    @("210×7=1470", "847×5=4235"),
    @("938×5=4690", "893×5=4465"),
    @("527×3=1581", "861×6=5166"),
    @("693×4=2772", "654×5=3270"),
    @("276×8=2208", "315×8=2520"),
    @("187×2=374",  "741×2=1482"),
    @("336×9=3024", "528×5=2640"),
    @("610×6=3660", "768×8=6144"),
    @("523×4=2092", "972×2=1944"),
    @("379×3=1137", "387×6=2322"),
    @("712×9=6408", "160×8=1280"),
    @("596×8=4768", "367×4=1468")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    Write-Host ("Replaced '" + $old + "' -> '" + $new + "': " + $found)
}
